$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Week 5 (rows 11-12) fully completed: mark remaining checkboxes as done
$ws.Range("G11").Value = $true
$ws.Range("J11").Value = $true
$ws.Range("G12").Value = $true

# Week 6 (row 13) getting started: mark first two checkboxes as done
$ws.Range("C13").Value = $true
$ws.Range("E13").Value = $true

# Update the active selection to reflect where the editor ended up
$ws.Range("P14").Select()
